# Update workbook for "Add data for 2022-03-23" commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-15"

# Update the label in column A for the March row to reflect new "through" date
$ws.Range("A4").Value = "March (through 03-15)"

# Update March row (row 4) values
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 19
$ws.Range("D4").Value = 30
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 16
$ws.Range("G4").Value = 27
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 66

# Update Total row (row 5) values
$ws.Range("B5").Value = 50
$ws.Range("C5").Value = 106
$ws.Range("D5").Value = 161
$ws.Range("E5").Value = 167
$ws.Range("F5").Value = 95
$ws.Range("G5").Value = 168
$ws.Range("H5").Value = 382
$ws.Range("I5").Value = 366
